# BI-1613 - Updating test files
# Insert a new "Term Type" column after "Full Name" (before "Description"),
# populate it for the data rows, and add a few blank trailing rows
# (matching the "missing categories" ontology template updates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at C (pushes Description.. over to the right)
$ws.Columns("C").Insert()

# 2) Header for the new column
$ws.Range("C1").Value = "Term Type"

# 3) Populate the new column's data rows (row 2 stays blank, matching source)
$ws.Range("C3").Value = "Phenotype"
$ws.Range("C4").Value = "Germplasm Attribute"
$ws.Range("C5").Value = "phenotype"
$ws.Range("C6").Value = "germplasm passport"

# 4) Give the new column's data cells (rows 2-6) the same bordered look as
#    the rest of the table body.
$bodyRange = $ws.Range("C2:C6")
$bodyRange.Borders.Color = 13553360
$bodyRange.Borders.LineStyle = 1
$bodyRange.Borders.Weight = 2

# 5) A few trailing blank rows were added underneath the table in column C
$tailRange = $ws.Range("C7:C9")
$tailRange.Interior.Pattern = -4142

# 6) Restore the view: scrolled back to the top-left, selection on B5
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B5").Select()
